$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 7730
$ws.Range("K3").Value = 7988
$ws.Range("D4").Value = 1954
$ws.Range("K4").Value = 1680
$ws.Range("K5").Value = 572
$ws.Range("K6").Value = 8913
$ws.Range("D7").Value = 27754
$ws.Range("K7").Value = 26883

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 490
$ws.Range("K3").Value = 531
$ws.Range("K4").Value = 100
$ws.Range("K6").Value = 591
$ws.Range("K7").Value = 1759

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K2").Value = 198
$ws.Range("K5").Value = 14
$ws.Range("K6").Value = 135
$ws.Range("K7").Value = 569

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K3").Value = 398
$ws.Range("K6").Value = 361
$ws.Range("K7").Value = 1129

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K3").Value = 153
$ws.Range("K7").Value = 447

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 251
$ws.Range("K3").Value = 293
$ws.Range("K6").Value = 268
$ws.Range("K7").Value = 884

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K3").Value = 148
$ws.Range("K7").Value = 625

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K2").Value = 120
$ws.Range("K7").Value = 451

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("K2").Value = 29
$ws.Range("K6").Value = 42
$ws.Range("K7").Value = 109

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K7").Value = 800
$ws.Range("K8").Value = 1759
$ws.Range("K11").Value = 472
$ws.Range("K19").Value = 772
$ws.Range("K20").Value = 663
$ws.Range("K27").Value = 259
$ws.Range("K29").Value = 1486
$ws.Range("K30").Value = 109
$ws.Range("K31").Value = 320
$ws.Range("K33").Value = 1129
$ws.Range("K36").Value = 348
$ws.Range("K37").Value = 884
$ws.Range("K42").Value = 996
$ws.Range("K44").Value = 216
$ws.Range("K47").Value = 184
$ws.Range("K48").Value = 335
$ws.Range("K52").Value = 693
$ws.Range("K60").Value = 161
$ws.Range("K61").Value = 22
$ws.Range("K62").Value = 9
$ws.Range("D63").Value = 352
$ws.Range("K63").Value = 74
$ws.Range("K65").Value = 625
$ws.Range("K67").Value = 1043
$ws.Range("K68").Value = 69
$ws.Range("K72").Value = 127
$ws.Range("K77").Value = 176
$ws.Range("K78").Value = 335
$ws.Range("K79").Value = 658
$ws.Range("K83").Value = 569
$ws.Range("K85").Value = 1239
$ws.Range("K88").Value = 286
$ws.Range("K89").Value = 403
$ws.Range("K90").Value = 257
$ws.Range("K92").Value = 99
$ws.Range("K93").Value = 106
$ws.Range("K94").Value = 363
$ws.Range("K95").Value = 447
$ws.Range("K97").Value = 219
$ws.Range("K98").Value = 144
$ws.Range("K99").Value = 451
$ws.Range("D101").Value = 27754
$ws.Range("K101").Value = 26883

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K6").Value = 126
$ws.Range("K7").Value = 320

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K3").Value = 379
$ws.Range("K7").Value = 1043

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("K3").Value = 32
$ws.Range("K4").Value = 14

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 417
$ws.Range("K3").Value = 523
$ws.Range("K6").Value = 441
$ws.Range("K7").Value = 1486

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K6").Value = 152
$ws.Range("K7").Value = 335

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K2").Value = 224
$ws.Range("K7").Value = 772

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("K4").Value = 12
$ws.Range("K7").Value = 216

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 266
$ws.Range("K3").Value = 291
$ws.Range("K6").Value = 380
$ws.Range("K7").Value = 996

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K2").Value = 102
$ws.Range("K7").Value = 335

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K6").Value = 172
$ws.Range("K7").Value = 658

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K2").Value = 226
$ws.Range("K3").Value = 206
$ws.Range("K6").Value = 191
$ws.Range("K7").Value = 663

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K3").Value = 111
$ws.Range("K7").Value = 348

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("K3").Value = 24
$ws.Range("K7").Value = 106

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K2").Value = 263
$ws.Range("K3").Value = 253
$ws.Range("K4").Value = 31
$ws.Range("K7").Value = 800

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K6").Value = 168
$ws.Range("K7").Value = 363

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("K4").Value = 18
$ws.Range("K6").Value = 57
$ws.Range("K7").Value = 184

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("K6").Value = 80
$ws.Range("K7").Value = 144

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K6").Value = 163
$ws.Range("K7").Value = 472

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("K6").Value = 122
$ws.Range("K7").Value = 219

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("K2").Value = 32
$ws.Range("K7").Value = 99

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K6").Value = 111
$ws.Range("K7").Value = 286

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K3").Value = 124
$ws.Range("K7").Value = 403

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K3").Value = 62
$ws.Range("K7").Value = 259

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K2").Value = 94
$ws.Range("K7").Value = 257

$ws = $wb.Worksheets.Item('North Park')
$ws.Range("K2").Value = 27
$ws.Range("K7").Value = 69

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("K3").Value = 46
$ws.Range("K7").Value = 161

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 409
$ws.Range("K3").Value = 428
$ws.Range("K5").Value = 36
$ws.Range("K6").Value = 305
$ws.Range("K7").Value = 1239

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("K6").Value = 59
$ws.Range("K7").Value = 127

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("K6").Value = 26
$ws.Range("K7").Value = 176

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K6").Value = 248
$ws.Range("K7").Value = 693

$ws = $wb.Worksheets.Item('Mount Greenwood')
$ws.Range("K3").Value = 7
$ws.Range("K7").Value = 22

$ws = $wb.Worksheets.Item('Millenium Park')
$ws.Range("K2").Value = 3

$ws = $wb.Worksheets.Item('Museum Campus')
$ws.Range("K7").Value = 9
